$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the "email" column (C) by shifting D:G left into C:F ---
# (Using a direct Columns.Delete() on column C corrupts the <cols> custom-width
#  entry in this runtime, so instead we copy the data+formats one column to the
#  left and then drop the now-redundant trailing column.)
$ws.Range("D1:G5").Copy()
$ws.Range("C1").PasteSpecial()
$ws.Columns.Item(7).Delete()

# --- Step 2: update the sample data rows (user 1-4 -> user 5-8, emails -> phone numbers) ---

# Row 2
$ws.Range("A2").Value = "testuser5"
$ws.Range("B2").Value = "lastname5"
$ws.Range("C2").Value = "'1234567898"
$ws.Range("D2").Value = "Selenium@123"
$ws.Range("E2").Value = "Student"
$ws.Range("F2").Value = "male"

# Row 3
$ws.Range("A3").Value = "testuser6"
$ws.Range("B3").Value = "lastname6"
$ws.Range("C3").Value = "'2234567898"
$ws.Range("D3").Value = "Selenium@124"
$ws.Range("E3").Value = "Doctor"
$ws.Range("F3").Value = "female"

# Row 4
$ws.Range("A4").Value = "testuser7"
$ws.Range("B4").Value = "lastname7"
$ws.Range("C4").Value = "'3123456789"
$ws.Range("D4").Value = "Selenium@125"
$ws.Range("E4").Value = "Engineer"
$ws.Range("F4").Value = "male"

# Row 5
$ws.Range("A5").Value = "testuser8"
$ws.Range("B5").Value = "lastname8"
$ws.Range("C5").Value = "'4234567898"
$ws.Range("D5").Value = "Selenium@126"
$ws.Range("E5").Value = "Scientist"
$ws.Range("F5").Value = "female"
